$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.342
$ws.Range("B3").Value = 5.828000000000001
$ws.Range("E3").Value = 16.3
$ws.Range("E12").Value = 17.845
$ws.Range("B14").Value = 5.857
$ws.Range("B21").Value = 9.873000000000001
$ws.Range("B23").Value = 7.728999999999999
$ws.Range("E24").Value = 16.585
$ws.Range("B25").Value = 5.103
$ws.Range("C25").Value = -11.52
$ws.Range("E25").Value = 17.288
$ws.Range("B26").Value = 6.181999999999999
$ws.Range("C27").Value = -13.466
$ws.Range("B29").Value = 5.656000000000001
$ws.Range("C31").Value = -13.163
$ws.Range("C39").Value = -12.858
$ws.Range("C48").Value = -11.332
$ws.Range("E50").Value = 16.331
$ws.Range("C51").Value = -11.266
$ws.Range("C52").Value = -11.303
$ws.Range("B53").Value = 6.909999999999999
$ws.Range("E53").Value = 17.343
$ws.Range("C55").Value = -13.106
$ws.Range("C56").Value = -13.537
$ws.Range("B57").Value = 4.864
$ws.Range("C57").Value = -13.852
$ws.Range("E57").Value = 16.423
$ws.Range("B59").Value = 4.741
$ws.Range("E61").Value = 16.796
$ws.Range("E63").Value = 17.562
$ws.Range("B69").Value = 5.354000000000001
$ws.Range("E70").Value = 17.77
$ws.Range("C73").Value = -12.752
$ws.Range("B79").Value = 5.6
$ws.Range("B83").Value = 5.702
$ws.Range("E86").Value = 16.46
$ws.Range("C89").Value = -11.627
$ws.Range("C90").Value = -12.91
$ws.Range("B91").Value = 5.600000000000001
$ws.Range("C92").Value = -10.878
$ws.Range("B93").Value = 5.673
$ws.Range("E98").Value = 16.18
$ws.Range("E100").Value = 16.879
$ws.Range("E102").Value = 16.318
